# Update cryptos list prices and 1h volume percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.260.67"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.803.39"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.93%  "

$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "617.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.07"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.801.59"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.00%  "

$ws.Range("E8").Value = "  +0.20%  "

$ws.Range("E9").Value = "  -0.07%  "

$ws.Range("E10").Value = "  +3.42%  "

$ws.Range("E11").Value = "  -3.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.494"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.00"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000257"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.444.53"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.799.41"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.239.81"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("E18").Value = "  -0.29%  "

$ws.Range("E19").Value = "  +0.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "515.21"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.66"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.59"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.71%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.729"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.52"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.23"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.39%  "

$ws.Range("E26").Value = "  -2.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.12"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000137"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +25.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("E30").Value = "  -3.10%  "

$ws.Range("E31").Value = "  +2.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.81"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.88%  "

$ws.Range("E33").Value = "  -1.03%  "

$ws.Range("E34").Value = "  -2.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.42%  "

$ws.Range("E36").Value = "  +0.52%  "

$ws.Range("E37").Value = "  +2.77%  "

$ws.Range("E38").Value = "  +0.87%  "

$ws.Range("E39").Value = "  +1.91%  "

$ws.Range("E40").Value = "  +2.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.42"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "44.11"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.77"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "422.14"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.066.58"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.75"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0366"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.59"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "136.22"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.48"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.23%  "
